# Electricity dispatch logit exponent calibration update
$wb = $excel.ActiveWorkbook

# Update the calibrated exponent value on the EDLE sheet (was -3, now 1)
$wsEDLE = $wb.Worksheets.Item("EDLE")
$wsEDLE.Range("B2").Value = 1

# Make the EDLE sheet the active/selected tab (previously "About" was selected)
$wsEDLE.Activate()
